$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new "Team" column header in AF1, matching the style of the other
# header cells (bold, centered, bordered) by copying the format from AE1.
$ws.Range("AF1").Value = "Team"
$ws.Range("AE1").Copy()
$ws.Range("AF1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Fill AF2:AF126 with the coach's club for every season row in this sheet.
$ws.Range("AF2:AF126").Value = "Brentford"
